$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.404.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4720"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2867"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06494"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.84"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "100.53"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07807"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.53"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7289"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.167"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.388.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.117.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.327"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.331"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.042"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.03"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.33%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.895"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09678"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.492"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.52%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04815"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6909"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.742"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.03"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.303"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.952"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4214"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8254"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.92"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.793"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.016"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05759"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "882.37"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.81%  "
